# Auto-applies the cryptos list update described in the commit diff.
# (Updated coin prices / volume percentages; rows 28-29 swapped content.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.697.41'
$ws.Range('E2').Value = '  +1.23%  '
$ws.Range('D3').Value = '1.644.90'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  +0.09%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '213.17'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.13%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.532'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  -2.40%  '
$ws.Range('E9').Value = '  -1.30%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.0613'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -0.20%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0890'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +1.61%  '
$ws.Range('D12').Value = '1.878.51'
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').Value = '1.647.67'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('E14').Value = '  -1.06%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.561'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -1.83%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '64.18'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -2.44%  '
$ws.Range('D17').Value = '27.706.12'
$ws.Range('E17').Value = '  +1.26%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '230.13'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('E19').Value = '  -0.18%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '7.67'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +2.83%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  -1.43%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '10.00'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +6.85%  '
$ws.Range('E24').Value = '  -3.58%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '149.15'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +1.51%  '
$ws.Range('B28').Value = 'BinanceUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '15.66'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('E31').Value = '  -2.60%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '3.30'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -0.06%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '3.17'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +1.04%  '
$ws.Range('D34').Value = '1.442.19'
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('E35').Value = '  +1.39%  '
$ws.Range('E36').Value = '  -1.12%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.570'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -0.12%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.883'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -2.85%  '
$ws.Range('E39').Value = '  -1.24%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.916'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +16.19%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '1.04'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('E42').Value = '  +0.10%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '5.67'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +2.15%  '
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('E45').Value = '  +1.70%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '65.59'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('D47').Value = '1.787.70'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('E48').Value = '  -1.23%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '86.23'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -2.13%  '
$ws.Range('E50').Value = '  -2.27%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '7.72'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -0.36%  '
